$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BAC-10")

# Row 11: ProveedorMonturas -> Proveedor (definition unchanged)
$ws.Range("A11").Value = "Proveedor"

# Row 10: ProveedorLentes -> Laboratorio ; definition updated
$ws.Range("A10").Value = "Laboratorio"
$ws.Range("B10").Value = "Empresa a la que se le compran los lentes formulados o con modificaciones y se le encargan los arreglos."

$ws.Range("D18").Select()
